$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 ("Rejection reasons"): a new value "coco" was typed into D21, which
# pushed the pre-existing values in columns E:AI one column to the right.
# Apply the resulting final values directly to each affected cell.

$ws.Range("D21").Value = "coco"
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = "'234"
$ws.Range("G21").Value = "'23"
$ws.Range("I21").Value = ""
$ws.Range("J21").Value = "insufficient sample size, poor methodology"
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = "not relevant, poor methodology"
$ws.Range("S21").Value = ""
$ws.Range("T21").Value = "poor methodology, insufficient sample size"
$ws.Range("U21").Value = ""
$ws.Range("V21").Value = "not relevant"
$ws.Range("X21").Value = ""
$ws.Range("Y21").Value = "not relevant, poor methodology"
$ws.Range("Z21").Value = "not relevant"
$ws.Range("AB21").Value = ""
$ws.Range("AC21").Value = "insufficient sample size, poor methodology"
$ws.Range("AH21").Value = ""
$ws.Range("AI21").Value = "poor methodology"
